# Add data for 2022-03-07
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet and update the header label to reflect the new "through" date
$ws.Name = "Through 2022-02-27"
$ws.Range("A3").Value = "February (through 02-27)"

# Update February row (row 3) with newly reported counts
$ws.Range("B3").Value = 10
$ws.Range("C3").Value = 36
$ws.Range("E3").Value = 51
$ws.Range("G3").Value = 67
$ws.Range("H3").Value = 120
$ws.Range("I3").Value = 136

# Update Total row (row 4) = January (row 2) + February (row 3)
$ws.Range("B4").Value = 36
$ws.Range("C4").Value = 87
$ws.Range("E4").Value = 137
$ws.Range("G4").Value = 141
$ws.Range("H4").Value = 337
$ws.Range("I4").Value = 295
